$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price (D) column target cells so numeric-looking
# strings (e.g. "29.209.43", "0.9994") are preserved exactly as text,
# matching the source data which stores these as inline strings.
$priceCells = @("D2", "D3", "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D51")
foreach ($cellAddr in $priceCells) {
    $ws.Range($cellAddr).NumberFormat = "@"
}

# Apply updated values
$ws.Range("D2").Value = '29.209.43'
$ws.Range("E2").Value = '  -1.00%  '
$ws.Range("D3").Value = '1.865.56'
$ws.Range("E3").Value = '  -0.70%  '
$ws.Range("D4").Value = '0.9994'
$ws.Range("D5").Value = '0.7099'
$ws.Range("E5").Value = '  -0.79%  '
$ws.Range("D6").Value = '241.57'
$ws.Range("E6").Value = '  -0.12%  '
$ws.Range("D7").Value = '0.9999'
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").Value = '0.3113'
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("D9").Value = '0.07655'
$ws.Range("E9").Value = '  -3.86%  '
$ws.Range("D10").Value = '24.67'
$ws.Range("E10").Value = '  -2.62%  '
$ws.Range("D11").Value = '0.08371'
$ws.Range("E11").Value = '  +1.01%  '
$ws.Range("D12").Value = '1.865.53'
$ws.Range("E12").Value = '  -0.48%  '
$ws.Range("D13").Value = '5.221'
$ws.Range("E13").Value = '  -1.28%  '
$ws.Range("D14").Value = '0.7104'
$ws.Range("E14").Value = '  -2.73%  '
$ws.Range("D15").Value = '91.36'
$ws.Range("E15").Value = '  +0.13%  '
$ws.Range("D16").Value = '29.228.07'
$ws.Range("E16").Value = '  -0.89%  '
$ws.Range("D17").Value = '5.944'
$ws.Range("E17").Value = '  +0.13%  '
$ws.Range("D18").Value = '243.65'
$ws.Range("E18").Value = '  -0.92%  '
$ws.Range("D19").Value = '0.000007822'
$ws.Range("E19").Value = '  -0.84%  '
$ws.Range("D20").Value = '2.114.47'
$ws.Range("E20").Value = '  +0.16%  '
$ws.Range("E21").Value = '  -2.02%  '
$ws.Range("D22").Value = '0.9993'
$ws.Range("E22").Value = '  -0.10%  '
$ws.Range("D23").Value = '7.873'
$ws.Range("E23").Value = '  -1.06%  '
$ws.Range("D24").Value = '0.9995'
$ws.Range("E24").Value = '  -0.10%  '
$ws.Range("D25").Value = '0.1638'
$ws.Range("E25").Value = '  +1.32%  '
$ws.Range("D26").Value = '163.15'
$ws.Range("E26").Value = '  -0.53%  '
$ws.Range("D27").Value = '8.951'
$ws.Range("E27").Value = '  -1.30%  '
$ws.Range("D28").Value = '18.50'
$ws.Range("E28").Value = '  +0.88%  '
$ws.Range("E29").Value = '  +0.42%  '
$ws.Range("D30").Value = '1.312'
$ws.Range("E30").Value = '  -3.58%  '
$ws.Range("D31").Value = '4.396'
$ws.Range("E31").Value = '  +0.05%  '
$ws.Range("D32").Value = '4.250'
$ws.Range("E32").Value = '  +3.32%  '
$ws.Range("D33").Value = '0.05152'
$ws.Range("E33").Value = '  -2.35%  '
$ws.Range("D34").Value = '0.7932'
$ws.Range("E34").Value = '  +8.92%  '
$ws.Range("D35").Value = '1.911'
$ws.Range("E35").Value = '  -2.68%  '
$ws.Range("D36").Value = '1.166'
$ws.Range("E36").Value = '  -2.83%  '
$ws.Range("D37").Value = '2.684'
$ws.Range("E37").Value = '  +0.15%  '
$ws.Range("D38").Value = '0.01857'
$ws.Range("E38").Value = '  -0.71%  '
$ws.Range("D39").Value = '2.708'
$ws.Range("E39").Value = '  -0.38%  '
$ws.Range("D40").Value = '1.160.58'
$ws.Range("E40").Value = '  -5.42%  '
$ws.Range("D41").Value = '6.331'
$ws.Range("E41").Value = '  +3.21%  '
$ws.Range("D42").Value = '0.8977'
$ws.Range("E42").Value = '  -1.62%  '
$ws.Range("D43").Value = '73.21'
$ws.Range("E43").Value = '  -1.00%  '
$ws.Range("D44").Value = '0.9992'
$ws.Range("E44").Value = '  -0.13%  '
$ws.Range("D45").Value = '103.47'
$ws.Range("E45").Value = '  +1.15%  '
$ws.Range("D46").Value = '2.012.82'
$ws.Range("E46").Value = '  -0.05%  '
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").Value = '1.779'
$ws.Range("E47").Value = '  -1.53%  '
$ws.Range("B48").Value = 'Mantle'
$ws.Range("C48").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D48").Value = '0.5142'
$ws.Range("E48").Value = '  -2.83%  '
$ws.Range("D49").Value = '9.336'
$ws.Range("E49").Value = '  -0.07%  '
$ws.Range("E50").Value = '  -0.98%  '
$ws.Range("D51").Value = '0.4292'
$ws.Range("E51").Value = '  -0.87%  '
